$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "...with a volume value higher than the previous one..."
#   "a volume value higher than " -> "a " + bold("volume value") + " higher than "
#   (pure text already exists -> only formatting is applied, via Range.Font.Bold,
#   which splits only the minimum runs necessary without disturbing neighbours)
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("with a volume value higher than the", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

if ($rng1.Find.Found) {
    $matchStart = $rng1.Start
    # "with " = 5 chars, "a " = 2 chars -> "volume value" begins right after
    $boldStart = $matchStart + 5 + 2
    $boldLen = "volume value".Length
    $boldEnd = $boldStart + $boldLen

    $boldRange = $d.Range($boldStart, $boldEnd)
    if ($boldRange.Text -eq "volume value") {
        $boldRange.Font.Bold = 1
        $boldRange.Font.BoldBi = 1
    }
}

# ---------------------------------------------------------------------------
# Edit 2: "...>0.35% combined with V change." ->
#          "...>0.35% combined with V change; = INCREASING VOLUME as the price is moving up."
#   New characters must be inserted, then bold is applied afterwards as a pure
#   formatting operation (keeps the collateral run-merging caused by the text
#   insertion to the minimum unavoidable amount).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(" combined with V change", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

if ($rng2.Find.Found) {
    $insertPoint = $rng2.End
    $newText = "= INCREASING VOLUME as the price is moving up"

    # Insert the literal text right before the following "." run.
    $dotRange = $d.Range($insertPoint, $insertPoint + 1)
    if ($dotRange.Text -eq ".") {
        $dotRange.InsertBefore("; " + $newText)
    } else {
        $zero = $d.Range($insertPoint, $insertPoint)
        $zero.InsertAfter("; " + $newText)
    }

    $boldStart = $insertPoint + 2
    $boldEnd = $boldStart + $newText.Length
    $boldRange2 = $d.Range($boldStart, $boldEnd)
    if ($boldRange2.Text -eq $newText) {
        $boldRange2.Font.Bold = 1
        $boldRange2.Font.BoldBi = 1
    }
}
